# Update countries & provincias Spain
# Refreshes the COVID-19 country table ("Pais" sheet) with a newer data
# pull: a handful of per-country metrics changed, the "last updated" time
# label moved from 16:05 to 16:35, and two pairs of countries swapped rank
# (because their "Casos totales" changed) in the table, which is sorted by
# that column descending.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp label (row 1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 13 de Mayo de 2020 a las 16:35"

# Estados Unidos (row 4) - updated totals
$ws.Cells.Item(4, 2).Value = 1411148
$ws.Cells.Item(4, 3).Value = 2512
$ws.Cells.Item(4, 4).Value = 298643
$ws.Cells.Item(4, 5).Value = 1028941
$ws.Cells.Item(4, 7).Value = 139
$ws.Cells.Item(4, 8).Value = 83564

# Argentina (row 55) - updated totals
$ws.Cells.Item(55, 4).Value = 2266
$ws.Cells.Item(55, 5).Value = 3976
$ws.Cells.Item(55, 7).Value = 2
$ws.Cells.Item(55, 8).Value = 321

# Rows 95/96 swap rank: "Consejo Danes para los Refugiados" overtakes
# "Mayotte" after its update, so it now sits in row 95 and Mayotte drops
# to row 96 (its own figures are unchanged).
$ws.Cells.Item(95, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(95, 2).Value = 1169
$ws.Cells.Item(95, 3).Value = 67
$ws.Cells.Item(95, 4).Value = 148
$ws.Cells.Item(95, 5).Value = 971
$ws.Cells.Item(95, 6).Value = 0
$ws.Cells.Item(95, 7).Value = 6
$ws.Cells.Item(95, 8).Value = 50

$ws.Cells.Item(96, 1).Value = "Mayotte"
$ws.Cells.Item(96, 2).Value = 1143
$ws.Cells.Item(96, 3).Value = 48
$ws.Cells.Item(96, 4).Value = 627
$ws.Cells.Item(96, 5).Value = 502
$ws.Cells.Item(96, 6).Value = 6
$ws.Cells.Item(96, 7).Value = 2
$ws.Cells.Item(96, 8).Value = 14

# Rows 115/116/117: "Kenia" overtakes "Mali" and "Uruguay" after its
# update, moving up to row 115; Mali and Uruguay each drop one row
# (their own figures are unchanged).
$ws.Cells.Item(115, 1).Value = "Kenia"
$ws.Cells.Item(115, 2).Value = 737
$ws.Cells.Item(115, 3).Value = 22
$ws.Cells.Item(115, 4).Value = 281
$ws.Cells.Item(115, 5).Value = 416
$ws.Cells.Item(115, 6).Value = 1
$ws.Cells.Item(115, 7).Value = 4
$ws.Cells.Item(115, 8).Value = 40

$ws.Cells.Item(116, 1).Value = "Mali"
$ws.Cells.Item(116, 2).Value = 730
$ws.Cells.Item(116, 3).Value = 0
$ws.Cells.Item(116, 4).Value = 398
$ws.Cells.Item(116, 5).Value = 292
$ws.Cells.Item(116, 6).Value = 0
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 40

$ws.Cells.Item(117, 1).Value = "Uruguay"
$ws.Cells.Item(117, 2).Value = 715
$ws.Cells.Item(117, 3).Value = 0
$ws.Cells.Item(117, 4).Value = 532
$ws.Cells.Item(117, 5).Value = 166
$ws.Cells.Item(117, 6).Value = 8
$ws.Cells.Item(117, 7).Value = 0
$ws.Cells.Item(117, 8).Value = 19

# Yemen (row 171) - updated totals
$ws.Cells.Item(171, 2).Value = 70
$ws.Cells.Item(171, 3).Value = 5
$ws.Cells.Item(171, 5).Value = 57
$ws.Cells.Item(171, 7).Value = 2
$ws.Cells.Item(171, 8).Value = 12
